$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OctoberRaw")

$ws.Range("A1").Value = "Library"
$ws.Range("A2").Value = "Atchison Public Library"
$ws.Range("B2").Value = 4038
$ws.Range("C2").Value = 1788
$ws.Range("D2").Value = 5826
$ws.Range("A3").Value = "Baldwin City Public Library"
$ws.Range("B3").Value = 2575
$ws.Range("C3").Value = 541
$ws.Range("D3").Value = 3116
$ws.Range("A4").Value = "Basehor Community Library"
$ws.Range("B4").Value = 8095
$ws.Range("C4").Value = 1106
$ws.Range("D4").Value = 9201
$ws.Range("A5").Value = "Bern Community Library"
$ws.Range("B5").Value = 154
$ws.Range("C5").Value = 82
$ws.Range("D5").Value = 236
$ws.Range("A6").Value = "Bonner Springs City Library"
$ws.Range("B6").Value = 5312
$ws.Range("C6").Value = 1249
$ws.Range("D6").Value = 6561
$ws.Range("A7").Value = "Burlingame Community Library"
$ws.Range("B7").Value = 450
$ws.Range("C7").Value = 219
$ws.Range("D7").Value = 669
$ws.Range("A8").Value = "Carbondale City Library"
$ws.Range("B8").Value = 629
$ws.Range("C8").Value = 81
$ws.Range("D8").Value = 710
$ws.Range("A9").Value = "Centralia Community Library"
$ws.Range("B9").Value = 300
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 352
$ws.Range("A10").Value = "Corning City Library"
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 17
$ws.Range("A11").Value = "Digital Content"
$ws.Range("A12").Value = "Doniphan County Library - Elwood"
$ws.Range("B12").Value = 159
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 179
$ws.Range("A13").Value = "Doniphan County Library - Highland"
$ws.Range("B13").Value = 253
$ws.Range("C13").Value = 182
$ws.Range("D13").Value = 435
$ws.Range("A14").Value = "Doniphan County Library - Troy"
$ws.Range("B14").Value = 544
$ws.Range("C14").Value = 119
$ws.Range("D14").Value = 663
$ws.Range("A15").Value = "Doniphan County Library - Wathena"
$ws.Range("B15").Value = 421
$ws.Range("C15").Value = 85
$ws.Range("D15").Value = 506
$ws.Range("A16").Value = "Effingham Community Library"
$ws.Range("B16").Value = 214
$ws.Range("C16").Value = 21
$ws.Range("D16").Value = 235
$ws.Range("A17").Value = "Eudora Community Library"
$ws.Range("B17").Value = 1694
$ws.Range("C17").Value = 769
$ws.Range("D17").Value = 2463
$ws.Range("A18").Value = "Everest, Barnes Reading Room"
$ws.Range("B18").Value = 121
$ws.Range("C18").Value = 74
$ws.Range("D18").Value = 195
$ws.Range("A19").Value = "Hiawatha, Morrill Public Library"
$ws.Range("B19").Value = 1515
$ws.Range("C19").Value = 619
$ws.Range("D19").Value = 2134
$ws.Range("A20").Value = "Highland Community College"
$ws.Range("B20").Value = 65
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 76
$ws.Range("A21").Value = "Holton, Beck-Bookman Library"
$ws.Range("B21").Value = 1803
$ws.Range("C21").Value = 441
$ws.Range("D21").Value = 2244
$ws.Range("A22").Value = "Horton Public Library"
$ws.Range("B22").Value = 128
$ws.Range("C22").Value = 48
$ws.Range("D22").Value = 176
$ws.Range("A23").Value = "Lansing Community Library"
$ws.Range("B23").Value = 1681
$ws.Range("C23").Value = 776
$ws.Range("D23").Value = 2457
$ws.Range("A24").Value = "Leavenworth Public Library"
$ws.Range("B24").Value = 8601
$ws.Range("C24").Value = 1710
$ws.Range("D24").Value = 10311
$ws.Range("A25").Value = "Linwood Community Library"
$ws.Range("B25").Value = 682
$ws.Range("C25").Value = 188
$ws.Range("D25").Value = 870
$ws.Range("A26").Value = "Louisburg Library"
$ws.Range("A27").Value = "Lyndon Carnegie Library"
$ws.Range("B27").Value = 450
$ws.Range("C27").Value = 234
$ws.Range("D27").Value = 684
$ws.Range("A28").Value = "McLouth Public Library"
$ws.Range("B28").Value = 149
$ws.Range("C28").Value = 41
$ws.Range("D28").Value = 190
$ws.Range("A29").Value = "Meriden-Ozawkie Public Library"
$ws.Range("B29").Value = 1255
$ws.Range("C29").Value = 583
$ws.Range("D29").Value = 1838
$ws.Range("A30").Value = "Northeast Kansas Library System"
$ws.Range("B30").Value = 14
$ws.Range("C30").Value = 36
$ws.Range("D30").Value = 50
$ws.Range("A31").Value = "Nortonville Public Library"
$ws.Range("B31").Value = 276
$ws.Range("C31").Value = 116
$ws.Range("D31").Value = 392
$ws.Range("A32").Value = "Osage City Library"
$ws.Range("B32").Value = 1590
$ws.Range("C32").Value = 485
$ws.Range("D32").Value = 2075
$ws.Range("A33").Value = "Osawatomie Public Library"
$ws.Range("B33").Value = 1036
$ws.Range("C33").Value = 428
$ws.Range("D33").Value = 1464
$ws.Range("A34").Value = "Oskaloosa Public Library"
$ws.Range("B34").Value = 521
$ws.Range("C34").Value = 181
$ws.Range("D34").Value = 702
$ws.Range("A35").Value = "Ottawa Library"
$ws.Range("B35").Value = 6101
$ws.Range("C35").Value = 912
$ws.Range("D35").Value = 7013
$ws.Range("A36").Value = "Overbrook Public Library"
$ws.Range("B36").Value = 945
$ws.Range("C36").Value = 198
$ws.Range("D36").Value = 1143
$ws.Range("A37").Value = "Paola Free Library"
$ws.Range("B37").Value = 3111
$ws.Range("C37").Value = 494
$ws.Range("D37").Value = 3605
$ws.Range("A38").Value = "Perry-Lecompton Community Library"
$ws.Range("B38").Value = 93
$ws.Range("C38").Value = 20
$ws.Range("D38").Value = 113
$ws.Range("A39").Value = "Pomona Community Library"
$ws.Range("B39").Value = 159
$ws.Range("C39").Value = 93
$ws.Range("D39").Value = 252
$ws.Range("A40").Value = "Prairie Hills Schools - Axtell Public School"
$ws.Range("B40").Value = 553
$ws.Range("C40").Value = 8
$ws.Range("D40").Value = 561
$ws.Range("A41").Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Range("B41").Value = 2340
$ws.Range("C41").Value = 110
$ws.Range("D41").Value = 2450
$ws.Range("A42").Value = "Prairie Hills Schools - Sabetha High School"
$ws.Range("B42").Value = 39
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 42
$ws.Range("A43").Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Range("B43").Value = 167
$ws.Range("C43").Value = 26
$ws.Range("D43").Value = 193
$ws.Range("A44").Value = "Prairie Hills Schools - Wetmore Academic Center (Permanently closed)"
$ws.Range("A45").Value = "Richmond Public Library"
$ws.Range("B45").Value = 373
$ws.Range("C45").Value = 89
$ws.Range("D45").Value = 462
$ws.Range("A46").Value = "Rossville Community Library"
$ws.Range("B46").Value = 1343
$ws.Range("C46").Value = 453
$ws.Range("D46").Value = 1796
$ws.Range("A47").Value = "Sabetha, Mary Cotton Library"
$ws.Range("B47").Value = 2899
$ws.Range("C47").Value = 1006
$ws.Range("D47").Value = 3905
$ws.Range("A48").Value = "Seneca Free Library"
$ws.Range("B48").Value = 1517
$ws.Range("C48").Value = 332
$ws.Range("D48").Value = 1849
$ws.Range("A49").Value = "Silver Lake Library"
$ws.Range("B49").Value = 1101
$ws.Range("C49").Value = 415
$ws.Range("D49").Value = 1516
$ws.Range("A50").Value = "Tonganoxie Public Library"
$ws.Range("B50").Value = 3121
$ws.Range("C50").Value = 1015
$ws.Range("D50").Value = 4136
$ws.Range("A51").Value = "Valley Falls, Delaware Township Library"
$ws.Range("B51").Value = 375
$ws.Range("C51").Value = 214
$ws.Range("D51").Value = 589
$ws.Range("A52").Value = "Wellsville City Library"
$ws.Range("B52").Value = 1386
$ws.Range("C52").Value = 408
$ws.Range("D52").Value = 1794
$ws.Range("A53").Value = "Wetmore Public Library"
$ws.Range("B53").Value = 167
$ws.Range("C53").Value = 120
$ws.Range("D53").Value = 287
$ws.Range("A54").Value = "Williamsburg Community Library"
$ws.Range("B54").Value = 390
$ws.Range("C54").Value = 29
$ws.Range("D54").Value = 419
$ws.Range("A55").Value = "Winchester Public Library"
$ws.Range("B55").Value = 295
$ws.Range("C55").Value = 329
$ws.Range("D55").Value = 624

$octSheet = $wb.Worksheets.Item("October")
$octSheet.Activate()
$octSheet.Range("B2").Select()
